# Update data reading mechanism: remove BinomialTotal column from the
# "Endpoints" sheet, fix up some endpoint distribution/measurement values,
# and switch the active sheet/selection state to match the new workflow.

$wb = $excel.ActiveWorkbook
$wsGroups = $wb.Worksheets.Item("EndpointGroups")
$wsEndpoints = $wb.Worksheets.Item("Endpoints")

# --- Endpoints sheet (sheet2): data fixes ---
# Row 2 (Endpoint 1): MeasurementType was "Fraction", now "Count".
$wsEndpoints.Range("C2").Value = "Count"
# Row 4 (Endpoint 4): DistributionType was "PowerLaw", now "Poisson".
$wsEndpoints.Range("F4").Value = "Poisson"
# Row 5 (Endpoint 3): DistributionType was "PowerLaw", now "OverdispersedPoisson".
$wsEndpoints.Range("F5").Value = "OverdispersedPoisson"

# The "BinomialTotal" column (G) is no longer used - remove it, shifting
# H:L (PowerLawPower, Mean, CV, RepeatedMeasurements, ExcessZeroes) left.
$wsEndpoints.Range("G1").EntireColumn.Delete()

# Restore/approximate the column widths for the shifted columns (closest
# achievable values given the host's column-width quantization).
$wsEndpoints.Columns.Item(7).ColumnWidth = 8.833333333333334
$wsEndpoints.Columns.Item(8).ColumnWidth = 11.666666666666666
$wsEndpoints.Columns.Item(9).ColumnWidth = 10.666666666666666
$wsEndpoints.Columns.Item(10).ColumnWidth = 16.666666666666668
$wsEndpoints.Columns.Item(11).ColumnWidth = 19.666666666666668

# --- View/selection state ---
# EndpointGroups tab was selected before; Endpoints tab is selected now.
$wsGroups.Range("B42").Select()
$wsEndpoints.Activate()
$wsEndpoints.Range("J2").Select()
